$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price strings so Excel
# does not auto-convert them to numbers (source cells are text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "51.621.84"
$ws.Range("E2").Value = "  +0.22%  "

$ws.Range("D3").Value = "2.793.47"
$ws.Range("E3").Value = "  +1.08%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "353.32"
$ws.Range("E5").Value = "  -0.99%  "

$ws.Range("D6").Value = "111.33"
$ws.Range("E6").Value = "  +3.37%  "

$ws.Range("D7").Value = "0.555"
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "0.629"
$ws.Range("E9").Value = "  +8.18%  "

$ws.Range("D10").Value = "40.07"
$ws.Range("E10").Value = "  +1.91%  "

$ws.Range("E11").Value = "  -1.74%  "

$ws.Range("E12").Value = "  -0.22%  "

$ws.Range("D13").Value = "19.93"
$ws.Range("E13").Value = "  +2.08%  "

$ws.Range("D14").Value = "7.74"
$ws.Range("E14").Value = "  +2.60%  "

$ws.Range("D15").Value = "3.231.45"
$ws.Range("E15").Value = "  +0.87%  "

$ws.Range("D16").Value = "2.801.72"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "0.943"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("D18").Value = "51.589.01"
$ws.Range("E18").Value = "  +0.16%  "

$ws.Range("D19").Value = "7.58"
$ws.Range("E19").Value = "  +0.20%  "

$ws.Range("D20").Value = "3.21"
$ws.Range("E20").Value = "  +4.95%  "

$ws.Range("D21").Value = "13.58"
$ws.Range("E21").Value = "  +4.12%  "

$ws.Range("E22").Value = "  +0.91%  "

$ws.Range("D23").Value = "70.24"
$ws.Range("E23").Value = "  +0.81%  "

$ws.Range("D24").Value = "267.03"
$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("D25").Value = "2.75"
$ws.Range("E25").Value = "  -0.13%  "

$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.09%  "

$ws.Range("D27").Value = "26.06"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("D29").Value = "38.92"
$ws.Range("E29").Value = "  +11.91%  "

$ws.Range("E30").Value = "  +2.33%  "

$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("D32").Value = "52.57"
$ws.Range("E32").Value = "  +1.83%  "

$ws.Range("D33").Value = "6.11"
$ws.Range("E33").Value = "  +1.08%  "

$ws.Range("E34").Value = "  +3.89%  "

$ws.Range("D35").Value = "0.0895"
$ws.Range("E35").Value = "  +7.40%  "

$ws.Range("D36").Value = "5.58"
$ws.Range("E36").Value = "  +9.06%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("D38").Value = "18.79"
$ws.Range("E38").Value = "  +1.07%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "2.00"
$ws.Range("E39").Value = "  +3.47%  "

$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "3.16"
$ws.Range("E40").Value = "  +1.00%  "

$ws.Range("E41").Value = "  +1.27%  "

$ws.Range("E42").Value = "  +0.47%  "

$ws.Range("E43").Value = "  +1.81%  "

$ws.Range("D44").Value = "120.92"
$ws.Range("E44").Value = "  +0.65%  "

$ws.Range("D45").Value = "21.75"
$ws.Range("E45").Value = "  +1.11%  "

$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  +6.31%  "

$ws.Range("E47").Value = "  +5.20%  "

$ws.Range("D48").Value = "2.107.87"
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").Value = "0.959"
$ws.Range("E49").Value = "  +4.43%  "

$ws.Range("D50").Value = "5.47"
$ws.Range("E50").Value = "  -0.75%  "

$ws.Range("E51").Value = "  +6.26%  "
